$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would be misread as a pure number by Excel (e.g. "1.00" -> 1,
# "3.20" -> 3.2) need the column pre-set to Text format so the literal string is kept.
$textCells = @(
    "D5",
    "D6",
    "D7",
    "D8",
    "D10",
    "D11",
    "D13",
    "D15",
    "D17",
    "D20",
    "D23",
    "D25",
    "D27",
    "D30",
    "D31",
    "D34",
    "D36",
    "D39",
    "D40",
    "D45",
    "D46",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Updated coin price (D) / 1h volume change (E) figures pulled from the latest run.
$ws.Range("D2").Value = '42.896.15'
$ws.Range("E2").Value = '  -4.78%  '
$ws.Range("D3").Value = '2.222.35'
$ws.Range("E3").Value = '  -5.82%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '317.41'
$ws.Range("E5").Value = '  +1.96%  '
$ws.Range("D6").Value = '100.15'
$ws.Range("E6").Value = '  -6.92%  '
$ws.Range("D7").Value = '0.591'
$ws.Range("E7").Value = '  -6.16%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  -6.70%  '
$ws.Range("D10").Value = '37.31'
$ws.Range("E10").Value = '  -8.36%  '
$ws.Range("D11").Value = '54.13'
$ws.Range("E11").Value = '  -2.55%  '
$ws.Range("E12").Value = '  -9.12%  '
$ws.Range("D13").Value = '7.81'
$ws.Range("E13").Value = '  -6.99%  '
$ws.Range("E14").Value = '  -2.72%  '
$ws.Range("D15").Value = '0.867'
$ws.Range("E15").Value = '  -10.74%  '
$ws.Range("D16").Value = '2.560.34'
$ws.Range("E16").Value = '  -5.83%  '
$ws.Range("D17").Value = '14.31'
$ws.Range("E17").Value = '  -5.60%  '
$ws.Range("D18").Value = '2.220.83'
$ws.Range("E18").Value = '  -6.09%  '
$ws.Range("D19").Value = '42.798.10'
$ws.Range("E19").Value = '  -4.92%  '
$ws.Range("D20").Value = '15.09'
$ws.Range("E20").Value = '  +6.34%  '
$ws.Range("D21").Value = '0.0₃0966'
$ws.Range("E21").Value = '  -8.54%  '
$ws.Range("E22").Value = '  -10.03%  '
$ws.Range("D23").Value = '65.54'
$ws.Range("E23").Value = '  -9.94%  '
$ws.Range("E24").Value = '  -9.39%  '
$ws.Range("D25").Value = '236.66'
$ws.Range("E25").Value = '  -8.19%  '
$ws.Range("E26").Value = '  -6.48%  '
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("E28").Value = '  -8.37%  '
$ws.Range("E29").Value = '  -4.62%  '
$ws.Range("D30").Value = '6.43'
$ws.Range("E30").Value = '  -10.16%  '
$ws.Range("D31").Value = '0.0909'
$ws.Range("E31").Value = '  -5.70%  '
$ws.Range("E32").Value = '  -7.74%  '
$ws.Range("E33").Value = '  -7.53%  '
$ws.Range("D34").Value = '156.79'
$ws.Range("E34").Value = '  -6.54%  '
$ws.Range("E35").Value = '  -6.08%  '
$ws.Range("D36").Value = '3.20'
$ws.Range("E36").Value = '  +10.50%  '
$ws.Range("E37").Value = '  +13.84%  '
$ws.Range("E38").Value = '  -5.74%  '
$ws.Range("D39").Value = '3.97'
$ws.Range("E39").Value = '  +0.96%  '
$ws.Range("D40").Value = '4.47'
$ws.Range("E40").Value = '  -4.42%  '
$ws.Range("E41").Value = '  -9.09%  '
$ws.Range("E42").Value = '  -7.23%  '
$ws.Range("D43").Value = '1.954.66'
$ws.Range("E43").Value = '  +3.98%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").Value = '12.45'
$ws.Range("E45").Value = '  -2.47%  '
$ws.Range("D46").Value = '88.81'
$ws.Range("E46").Value = '  -10.90%  '
$ws.Range("E47").Value = '  -8.60%  '
$ws.Range("D48").Value = '5.38'
$ws.Range("E48").Value = '  -4.12%  '
$ws.Range("D49").Value = '76.41'
$ws.Range("E49").Value = '  -5.77%  '
$ws.Range("D50").Value = '60.64'
$ws.Range("E50").Value = '  -12.31%  '
$ws.Range("D51").Value = '0.879'
$ws.Range("E51").Value = '  +20.64%  '
